# Revised Statistics, Results, & first draft Discussion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Overall" row EQ5D mean / SE figures
$ws.Range("B3").Value = 81.75
$ws.Range("C3").Value = 7.8

# Restore the selection to C4 (as left by the author after editing)
$ws.Range("C4").Select()
